$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Ação"
$ws.Range("J1").Value = "Status"
$ws.Range("K1").Value = "ID REGISTRO"
$ws.Range("L1").Value = "ID SERVIÇO"

$ws.Range("I2").Value = "CADASTRAR"
$ws.Range("I3").Value = "CADASTRAR"
